# IRFC.NS.xlsx (1wk) update:
#  - Q59 (detect_structure) 2 -> 0
#  - R179 / R180 (backup) "" -> 0
#  - append 22 new weekly OHLCV rows (181-202), Jul-Nov 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cells ---
$ws.Cells.Item(59, 17).Value = 0   # Q59: detect_structure
$ws.Cells.Item(179, 18).Value = 0  # R179: backup
$ws.Cells.Item(180, 18).Value = 0  # R180: backup

# --- Append 22 new weekly rows (181-202), columns A:Q ---
# Datetime,Open,High,Low,Close,AdjClose,Volume,Year,Month,Day,Hour,Minute,Second,Week,isPivot,two_line_structure,detect_structure
$rowsText = @(
    "45474,173.7899932861328,191.5,170.0099945068359,188.2299957275391,186.4824981689453,340150588,2024,7,1,0,0,0,27,0,0,0",
    "45481,193.5,224.3999938964844,192.1999969482422,217.4299926757812,215.4114074707031,1007398123,2024,7,8,0,0,0,28,0,0,0",
    "45488,220.3999938964844,229,200.6999969482422,205.7700042724609,203.8596649169922,472487894,2024,7,15,0,0,0,29,1,0,0",
    "45495,205.7700042724609,209.1999969482422,164.1499938964844,183.4400024414062,181.7369689941406,429584237,2024,7,22,0,0,0,30,2,0,0",
    "45502,185.9199981689453,198.3000030517578,184.5800018310547,192.3300018310547,190.54443359375,271833689,2024,7,29,0,0,0,31,0,0,0",
    "45509,180.0299987792969,187,175.1300048828125,179.7899932861328,178.120849609375,246917424,2024,8,5,0,0,0,32,0,0,0",
    "45516,180,189.4499969482422,176.6999969482422,178.7700042724609,177.1103363037109,158052103,2024,8,12,0,0,0,33,0,0,0",
    "45523,181.1000061035156,187.3300018310547,178.1999969482422,183.4100036621094,181.7072601318359,94419519,2024,8,19,0,0,0,34,0,0,0",
    "45530,184,184.1999969482422,178,178.6199951171875,177.6505432128906,68712836,2024,8,26,0,0,0,35,0,0,2",
    "45537,178.8500061035156,180.25,168.1000061035156,169.8500061035156,168.9281616210938,75255427,2024,9,2,0,0,0,36,0,0,1",
    "45544,169.8399963378906,171.8399963378906,163.5200042724609,167.25,166.3422698974609,82067085,2024,9,9,0,0,0,37,0,0,0",
    "45551,167.25,167.4499969482422,151.6999969482422,161.4299926757812,160.5538482666016,149276895,2024,9,16,0,0,0,38,0,0,0",
    "45558,163,163,154.8300018310547,156.8099975585938,155.9589233398438,75330259,2024,9,23,0,0,0,39,0,0,0",
    "45565,156.6999969482422,159.3000030517578,146.8800048828125,152.0700073242188,151.2446594238281,82753825,2024,9,30,0,0,0,40,0,0,0",
    "45572,152.0700073242188,155.6900024414062,141.5099945068359,152.3999938964844,151.5728607177734,129918628,2024,10,7,0,0,0,41,0,0,0",
    "45579,153.5,154.5,145.3200073242188,147.6799926757812,146.8784637451172,62939825,2024,10,14,0,0,0,42,0,0,0",
    "45586,148,149,132.8000030517578,134.4299926757812,133.7003784179688,89065858,2024,10,21,0,0,0,43,2,0,0",
    "45593,134.4299926757812,161,134.0800018310547,158.1100006103516,157.2518768310547,209539983,2024,10,28,0,0,0,44,0,0,0",
    "45600,159,159.2400054931641,147,148.5399932861328,147.7338104248047,112584092,2024,11,4,0,0,0,45,0,0,0",
    "45607,147.5099945068359,151.1999969482422,139.3000030517578,139.6600036621094,138.9020080566406,63185760,2024,11,11,0,0,0,46,0,0,0",
    "45614,140.9799957275391,146.9600067138672,137.8000030517578,142.1499938964844,142.1499938964844,84756760,2024,11,18,0,0,0,47,0,0,0",
    "45621,149,155.6499938964844,145.0500030517578,149.3399963378906,149.3399963378906,103770647,2024,11,25,0,0,0,48,0,0,0"
)

$startRow = 181
$nRows = $rowsText.Count
$nCols = 17
$newData = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    $parts = $rowsText[$i].Split(",")
    for ($j = 0; $j -lt $nCols; $j++) {
        $newData[$i, $j] = [double]$parts[$j]
    }
}

$endRow = $startRow + $nRows - 1
$ws.Range("A$startRow`:Q$endRow").Value = $newData

# Column A carries the sheet's datetime display format (reuses the existing style)
$ws.Range("A$startRow`:A$endRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"
